$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New feedback rows appended below the existing data (rows 2-5 already present).
$rows = @(
    @("20BEE2001", "Merwin S",    "2020", "Ok",         4, "2026-02-05T04:52:39.104Z"),
    @("20BEE2001", "Merwin S",    "2020", "excellente", 5, "2026-02-05T04:57:00.418Z"),
    @("20BEE2001", "Merwin S",    "2020", "Ok",         5, "2026-02-05T05:03:57.480Z"),
    @("20BEE0001", "Merwin S",    "2020", "Ok",         5, "2026-02-05T05:11:22.267Z"),
    @("20BEE5984", "Ritwik",      "2020", "Not bad",    3, "2026-02-05T05:12:03.250Z"),
    @("20BEE2001", "Jershwin S",  "2020", "Niche",      5, "2026-02-05T05:16:10.354Z"),
    @("20BEE0001", "Merwin S",    "2020", "Manual",     5, "2026-02-05T05:18:25.149Z"),
    @("20BEE2001", "Merwin",      "2020", "Ok",         2, "5/2/2026, 11:12:29 am")
)

$startRow = 6
$lastRow = $startRow + $rows.Count - 1

# Column C ("DeptYear") holds purely-numeric text like "2020" - force the
# range to text format first so Excel doesn't silently coerce it to a number.
$ws.Range("C$startRow`:C$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
    $ws.Cells.Item($r, 6).Value = $data[5]
}
